$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new shared-string values in the exact order needed so the
# resulting sharedStrings.xml append-order lines up with the source file.
$ws.Range("F2").Value = "BAkso008.JPG"
$ws.Range("B5").Value = "Thangkas - Big Tsagil"
$ws.Range("B4").Value = "Malas - Arm Mala"
$ws.Range("F4").Value = "MAM020.JPG"
$ws.Range("B3").Value = "Buddhas - Shakyamuni"
$ws.Range("F3").Value = "BSha001.JPG, BSha002.JPG, BSha003.JPG"
$ws.Range("F5").Value = "TBT020.JPG, TBT019.JPG, TBT018.JPG"

# These reuse already-existing shared strings (order doesn't matter).
$ws.Range("C3").Value = "manjusri"
$ws.Range("C4").Value = "kleine mala"
$ws.Range("C5").Value = "große Thangka"

# Remove old rows 6 and 7 (their data is no longer needed).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Apply Text number format to column F (header cell keeps explicit style).
$ws.Range("F1").NumberFormat = "@"

# Restore selection to match the saved state.
$ws.Range("B7").Select()
